# Update cryptocurrency Price (D) and Volume(1h) (E) columns per the
# latest GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.852.64"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "2.499.50"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.56%  "

$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.37%  "

$ws.Range("D9").Value = "2.501.32"
$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("E10").Value = "  +2.13%  "

$ws.Range("E11").Value = "  -2.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.329"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("D14").Value = "2.946.58"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").Value = "58.767.92"
$ws.Range("E15").Value = "  +1.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").Value = "2.500.69"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.16%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.407"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.27%  "

$ws.Range("D30").Value = "0.0₃0753"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("E31").Value = "  +1.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("E33").Value = "  -0.91%  "

$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.22%  "

$ws.Range("E37").Value = "  -3.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("E39").Value = "  +3.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.821"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("E42").Value = "  +1.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "273.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "131.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0934"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.40%  "

$ws.Range("E49").Value = "  +1.72%  "

$ws.Range("E50").Value = "  -1.35%  "

$ws.Range("D51").Value = "1.749.05"
$ws.Range("E51").Value = "  +0.30%  "
